$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Update the status text for the 6be93617-... file from "Ready for handoff" to
# "Handback transform failed" everywhere it's referenced (Overview zh-cn/de-de
# status columns, and the Status column on each language sheet).
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff file name mismatch error detail for each locale.
$wsZhCn.Range("P3").Value = "Handback file name: 5lgmwrig.h4w is different with handoff file name: 6be93617-28b0-4d02-b90d-046e74e7ccb5.6f17e061bfcb46548cad038c875558de364e6813.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 5lgmwrig.h4w is different with handoff file name: 6be93617-28b0-4d02-b90d-046e74e7ccb5.6f17e061bfcb46548cad038c875558de364e6813.de-de."

# Widen the Error Detail column (column 16 / P) so the new message is readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
